# Generate Report for Handback
# Update the handoff/handback timestamps for the 6e0a8b13-... entry
# (row 3) on both the "zh-cn" and "de-de" worksheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-22 14:47:22"
$wsZhCn.Range("H3").Value = "2016-03-22 14:47:46"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-22 14:47:26"
$wsDeDe.Range("H3").Value = "2016-03-22 14:47:52"
